$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("store")

# Clear H and I columns on rows 5-7 (drop the trailing commentary columns)
$ws.Range("H5:I7").Clear()

# Row 8 (A8:C8 = 1,1,0, default style) moves up to become the new row 5
# content; copy it (values + formatting) onto A5:C5, then clear the rest.
$ws.Range("A8:C8").Copy($ws.Range("A5:C5"))

$ws.Range("D5:G5").ClearContents()
$ws.Range("A6:G7").ClearContents()

# Remove the old trailing row 8 entirely
$ws.Rows.Item(8).Delete()

# Update the active selection to match the new layout
$ws.Range("A5:C5").Select()
